# Apply crypto price/volume updates per the commit diff.
# Column D (Price) values are forced to text via an apostrophe-prefixed
# Formula assignment (then the style is reset to Normal) so that Excel
# does not reinterpret dotted/odd-format price strings as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Formula = "'" + '59.936.66'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.63%  '

$c = $ws.Range("D3")
$c.Formula = "'" + '2.534.74'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +1.91%  '

$ws.Range("E4").Value = '  -0.06%  '

$c = $ws.Range("D5")
$c.Formula = "'" + '543.59'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.70%  '

$c = $ws.Range("D6")
$c.Formula = "'" + '145.01'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.19%  '

$ws.Range("E7").Value = '  -0.31%  '

$c = $ws.Range("D8")
$c.Formula = "'" + '0.572'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.80%  '

$c = $ws.Range("D9")
$c.Formula = "'" + '2.568.65'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +2.12%  '

$c = $ws.Range("D10")
$c.Formula = "'" + '0.101'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.70%  '

$ws.Range("E11").Value = '  +1.09%  '

$c = $ws.Range("D12")
$c.Formula = "'" + '5.54'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +3.63%  '

$ws.Range("E13").Value = '  +1.25%  '

$c = $ws.Range("D14")
$c.Formula = "'" + '2.990.73'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.20%  '

$c = $ws.Range("D15")
$c.Formula = "'" + '23.88'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.02%  '

$c = $ws.Range("D16")
$c.Formula = "'" + '59.861.20'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.49%  '

$ws.Range("E17").Value = '  +2.27%  '

$c = $ws.Range("D18")
$c.Formula = "'" + '2.559.44'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.91%  '

$c = $ws.Range("D19")
$c.Formula = "'" + '11.34'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.52%  '

$ws.Range("E20").Value = '  -0.91%  '

$c = $ws.Range("D21")
$c.Formula = "'" + '327.98'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.30%  '

$ws.Range("E22").Value = '  -0.18%  '

$c = $ws.Range("D23")
$c.Formula = "'" + '5.93'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +2.30%  '

$c = $ws.Range("D24")
$c.Formula = "'" + '62.45'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.87%  '

$c = $ws.Range("D25")
$c.Formula = "'" + '0.441'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.32%  '

$c = $ws.Range("D26")
$c.Formula = "'" + '0.166'
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Formula = "'" + '0.993'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.69%  '

$c = $ws.Range("D28")
$c.Formula = "'" + '8.01'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +2.11%  '

$c = $ws.Range("D29")
$c.Formula = "'" + '7.10'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.68%  '

$c = $ws.Range("D30")
$c.Formula = "'" + '0.0' + [char]0x2083 + '0800'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.64%  '

$ws.Range("E31").Value = '  -0.60%  '

$c = $ws.Range("D32")
$c.Formula = "'" + '1.23'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -4.67%  '

$c = $ws.Range("D33")
$c.Formula = "'" + '163.19'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.47%  '

$ws.Range("E34").Value = '  +4.98%  '

$ws.Range("E35").Value = '  -0.12%  '

$c = $ws.Range("D36")
$c.Formula = "'" + '18.78'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.25%  '

$c = $ws.Range("D37")
$c.Formula = "'" + '4.47'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.80%  '

$c = $ws.Range("D38")
$c.Formula = "'" + '1.63'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.75%  '

$c = $ws.Range("D39")
$c.Formula = "'" + '5.71'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -4.46%  '

$c = $ws.Range("D40")
$c.Formula = "'" + '37.23'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +1.44%  '

$c = $ws.Range("D41")
$c.Formula = "'" + '302.27'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -4.19%  '

$c = $ws.Range("D42")
$c.Formula = "'" + '0.839'
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Formula = "'" + '3.74'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.70%  '

$ws.Range("E44").Value = '  +1.56%  '

$ws.Range("E45").Value = '  -0.10%  '

$c = $ws.Range("D46")
$c.Formula = "'" + '10.84'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.47%  '

$c = $ws.Range("D47")
$c.Formula = "'" + '19.08'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +2.63%  '

$c = $ws.Range("D48")
$c.Formula = "'" + '0.0938'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.17%  '

$c = $ws.Range("D49")
$c.Formula = "'" + '124.71'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.45%  '

$c = $ws.Range("D50")
$c.Formula = "'" + '0.0521'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.26%  '

$ws.Range("E51").Value = '  -1.00%  '
